$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C1: same text/style as B1 ("readProperties") ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = $ws.Range("B1").Value()

# --- C2: same (blank) style as B2/A2 ---
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# --- C3: same style as B3, then tweak alignment to left, then set the JSON text ---
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").HorizontalAlignment = -4131
$json = "{`n  ""target"": ""json"",`n  ""value"": ""src/test/resources/readproperties/readProperties2.json""`n}"
$ws.Range("C3").Value = $json

# --- column widths (B loses its auto bestFit, new column C gets a similar custom width) ---
$ws.Columns.Item(2).ColumnWidth = 41.857142857142854
$ws.Columns.Item(3).ColumnWidth = 42.142857142857146

# --- row heights ---
$ws.Rows.Item(1).RowHeight = 31.5
$ws.Rows.Item(3).RowHeight = 78.75

# --- selection moves to F5 ---
$ws.Range("F5").Select() | Out-Null
